$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new title row above the existing header row, pushing the
# --- existing header ("Code","Name","Unit","Price","Initial Stock","Stock") down to row 2.
$ws.Rows.Item(1).Insert()

# --- The inserted-above row leaves row 2 carrying the old row-1 formatting
# --- (bold/custom row format); clear it back to the sheet default before
# --- re-populating the header labels.
$ws.Rows.Item(2).ClearFormats()

# --- Row 2: rewrite the header labels - two new columns are introduced
# --- ("No" and "Product Category Code" before "Code", "Part Number" /
# --- "Brand" replacing "Initial Stock").
$ws.Range("A2").Value2 = "No"
$ws.Range("B2").Value2 = "Code"
$ws.Range("C2").Value2 = "Product Category Code"
$ws.Range("D2").Value2 = "Name"
$ws.Range("E2").Value2 = "Unit"
$ws.Range("F2").Value2 = "Price"
$ws.Range("G2").Value2 = "Part Number"
$ws.Range("H2").Value2 = "Brand"
$ws.Range("I2").Value2 = "Stock"

# --- Row 1: merged title cell spanning the whole table width, centered.
$ws.Range("A1").Value2 = "Template Product"
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").Merge()

# --- Column widths (character units converted from the authored pixel widths).
$ws.Columns.Item(1).ColumnWidth = 2.666666666666667
$ws.Columns.Item(2).ColumnWidth = 8.333333333333332
$ws.Columns.Item(3).ColumnWidth = 19.666666666666668
$ws.Columns.Item(4).ColumnWidth = 11.5
$ws.Columns.Item(5).ColumnWidth = 8.333333333333332
$ws.Columns.Item(6).ColumnWidth = 11.5
$ws.Columns.Item(7).ColumnWidth = 10.833333333333332
$ws.Columns.Item(8).ColumnWidth = 8.333333333333332

# --- Selection matches the authored diff.
$ws.Range("F9").Select()
